# Generate Report for Handoff
# Replaces the two tracked source files with a new pair:
#   4f94797c-8777-4592-8da8-a09997000153.md -> 62ad09a9-8825-480a-b9c6-c9b050f56804.md
#   d332a154-1705-4dfb-b938-0752c3ccafb3.md -> ffff3158594b-454b-4fba-b695-e0ab46dfeae7.md
# and updates status/timestamps/handoff-xliff bookkeeping to reflect a
# fresh "handoff" pass (vs. the previous "handback" state).

$wb = $excel.ActiveWorkbook

function Set-HyperlinkDisplay($ws, $targetAddr, $text) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $targetAddr) {
            $h.TextToDisplay = $text
            return
        }
    }
}

function Remove-HyperlinkAt($ws, $targetAddr) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $targetAddr) {
            $h.Delete()
            return
        }
    }
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "62ad09a9-8825-480a-b9c6-c9b050f56804.md"
$wsOverview.Range("B2").Value = "e2e\62ad09a9-8825-480a-b9c6-c9b050f56804.md"
$wsOverview.Range("C2").Value = ".md"
$wsOverview.Range("D2").Value = ""
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-28 23:01:41"

$wsOverview.Range("A3").Value = "ffff3158594b-454b-4fba-b695-e0ab46dfeae7.md"
$wsOverview.Range("B3").Value = "e2e\ffff3158594b-454b-4fba-b695-e0ab46dfeae7.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-28 23:01:41"

Set-HyperlinkDisplay $wsOverview '$B$2' "e2e\62ad09a9-8825-480a-b9c6-c9b050f56804.md"
Set-HyperlinkDisplay $wsOverview '$B$3' "e2e\ffff3158594b-454b-4fba-b695-e0ab46dfeae7.md"

$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "62ad09a9-8825-480a-b9c6-c9b050f56804.md"
$wsZhCn.Range("B2").Value = ".md"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("D2").Value = "e2e"
$wsZhCn.Range("E2").Value = "ht"
$wsZhCn.Range("F2").Value = "False"
$wsZhCn.Range("G2").Value = "62ad09a9-8825-480a-b9c6-c9b050f56804.5d3ee14a47ac78dc4ac2cae4ea69a21bef042b81.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-28 23:01:36"
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L2").Value = ""
$wsZhCn.Range("M2").Value = "True"
$wsZhCn.Range("N2").Value = ""
$wsZhCn.Range("O2").Value = "False"
$wsZhCn.Range("P2").Value = ""

$wsZhCn.Range("A3").Value = "ffff3158594b-454b-4fba-b695-e0ab46dfeae7.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "True"
$wsZhCn.Range("G3").Value = "62ad09a9-8825-480a-b9c6-c9b050f56804.5d3ee14a47ac78dc4ac2cae4ea69a21bef042b81.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-28 23:01:36"
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Range("P3").Value = ""

Set-HyperlinkDisplay $wsZhCn '$A$2' "62ad09a9-8825-480a-b9c6-c9b050f56804.md"
Set-HyperlinkDisplay $wsZhCn '$A$3' "ffff3158594b-454b-4fba-b695-e0ab46dfeae7.md"

Remove-HyperlinkAt $wsZhCn '$I$2'
Remove-HyperlinkAt $wsZhCn '$I$3'

$wsZhCn.Range("I2").Style = "Normal"
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("J2").Value = ""
$wsZhCn.Range("I3").Style = "Normal"
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""

$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsZhCn.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZhCn.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "62ad09a9-8825-480a-b9c6-c9b050f56804.md"
$wsDeDe.Range("B2").Value = ".md"
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("D2").Value = "e2e"
$wsDeDe.Range("E2").Value = "ht"
$wsDeDe.Range("F2").Value = "False"
$wsDeDe.Range("G2").Value = "62ad09a9-8825-480a-b9c6-c9b050f56804.5d3ee14a47ac78dc4ac2cae4ea69a21bef042b81.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-28 23:01:41"
$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L2").Value = ""
$wsDeDe.Range("M2").Value = "True"
$wsDeDe.Range("N2").Value = ""
$wsDeDe.Range("O2").Value = "False"
$wsDeDe.Range("P2").Value = ""

$wsDeDe.Range("A3").Value = "ffff3158594b-454b-4fba-b695-e0ab46dfeae7.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "True"
$wsDeDe.Range("G3").Value = "62ad09a9-8825-480a-b9c6-c9b050f56804.5d3ee14a47ac78dc4ac2cae4ea69a21bef042b81.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-28 23:01:41"
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""

Set-HyperlinkDisplay $wsDeDe '$A$2' "62ad09a9-8825-480a-b9c6-c9b050f56804.md"
Set-HyperlinkDisplay $wsDeDe '$A$3' "ffff3158594b-454b-4fba-b695-e0ab46dfeae7.md"

Remove-HyperlinkAt $wsDeDe '$I$2'
Remove-HyperlinkAt $wsDeDe '$I$3'

$wsDeDe.Range("I2").Style = "Normal"
$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("J2").Value = ""
$wsDeDe.Range("I3").Style = "Normal"
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""

$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsDeDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDeDe.Columns.Item(10).ColumnWidth = 21.7054770333426
